$p = $ppt.ActivePresentation
$oldUrl = "https://bndr.it/nwjmb"
$newUrl = "https://bndr.it/b5dn7"
$changed = 0

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }
        $tf = $shape.TextFrame
        if (-not $tf.HasText) { continue }
        $tr = $tf.TextRange
        for ($pi = 1; $pi -le 50; $pi++) {
            $para = $tr.Paragraphs($pi)
            $paraText = $para.Text
            if ([string]::IsNullOrEmpty($paraText)) { continue }
            if ($paraText -notlike "*$oldUrl*") { continue }
            for ($ri = 1; $ri -le 20; $ri++) {
                $run = $para.Runs($ri)
                $runText = $run.Text
                if ([string]::IsNullOrEmpty($runText)) { continue }
                if ($runText -eq $oldUrl) {
                    $run.Text = $newUrl
                    $changed++
                }
            }
        }
    }
}
Write-Output "changed = $changed"
